$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "total : 1,9Gi"

$ws.Range("E2").Value = "PID: 1579, COMMAND: gnome-shell, %MEM: 16.1`nPID: 1805, COMMAND: gnome-software, %MEM: 5.0`nPID: 1789, COMMAND: evolution-alarm, %MEM: 4.6`nPID: 2053, COMMAND: gsd-xsettings, %MEM: 4.0`nPID: 2028, COMMAND: Xwayland, %MEM: 3.3"

$ws.Range("F2").Value = "HOME : /home/rudy`nUSER : rudy`nPATH : /usr/local/bin:/usr/bin:/bin:/usr/games"

$ws.Range("G2").Value = "NAME: ├─sda1, SIZE: 19G, TYPE: part, MOUNTPOINT: /`nNAME: └─sda5, SIZE: 975M, TYPE: part, MOUNTPOINT: [SWAP]"

$ws.Range("H2").Value = "/dev : {'size': '934M', 'used': '0', 'avail': '934M', 'pcent': '0%'}`n/run : {'size': '194M', 'used': '1,4M', 'avail': '192M', 'pcent': '1%'}`n/ : {'size': '19G', 'used': '5,3G', 'avail': '13G', 'pcent': '30%'}`n/dev/shm : {'size': '967M', 'used': '0', 'avail': '967M', 'pcent': '0%'}`n/run/lock : {'size': '5,0M', 'used': '8,0K', 'avail': '5,0M', 'pcent': '1%'}`n/run/user/1000 : {'size': '194M', 'used': '96K', 'avail': '194M', 'pcent': '1%'}"
